# Auto-generated edit script applying the market-price refresh diff
# (static cached numeric values only -- no formulas in this workbook)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3415.577
$ws.Range("I64").Value = 3217.4783
$ws.Range("J64").Value = 4934.3335
$ws.Range("K64").Value = 3217.4783
$ws.Range("L64").Value = 4934.3335
$ws.Range("M64").Value = -2969.4783
$ws.Range("N64").Value = -5430.3335
$ws.Range("H67").Value = 3415.577
$ws.Range("I67").Value = 3217.4783
$ws.Range("J67").Value = 4934.3335
$ws.Range("K67").Value = 3217.4783
$ws.Range("L67").Value = 4934.3335
$ws.Range("M67").Value = -2359.4783
$ws.Range("N67").Value = -6650.3335
$ws.Range("H74").Value = 9469.625
$ws.Range("I74").Value = 11438.25
$ws.Range("J74").Value = 7501
$ws.Range("K74").Value = 11438.25
$ws.Range("L74").Value = 7501
$ws.Range("M74").Value = -10502.25
$ws.Range("N74").Value = -9373
$ws.Range("H77").Value = 9469.625
$ws.Range("I77").Value = 11438.25
$ws.Range("J77").Value = 7501
$ws.Range("K77").Value = 57191.25
$ws.Range("L77").Value = 37505
$ws.Range("M77").Value = -52511.25
$ws.Range("N77").Value = -46865
$ws.Range("H132").Value = 34246.613
$ws.Range("I132").Value = 36270.516
$ws.Range("K132").Value = 108811.548
$ws.Range("M132").Value = -106281.548
$ws.Range("H135").Value = 3236.5454
$ws.Range("I135").Value = 1556.1818
$ws.Range("J135").Value = 6597.273
$ws.Range("K135").Value = 14005.6362
$ws.Range("L135").Value = 59375.457
$ws.Range("M135").Value = -11470.6362
$ws.Range("N135").Value = -64445.457
$ws.Range("H141").Value = 5668.4614
$ws.Range("I141").Value = 6162.727
$ws.Range("J141").Value = 2950
$ws.Range("K141").Value = 18488.181
$ws.Range("L141").Value = 8850
$ws.Range("M141").Value = -13308.181
$ws.Range("N141").Value = -19210

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1720.174
$ws.Range("I61").Value = 1023
$ws.Range("J61").Value = 2804.6667
$ws.Range("K61").Value = 1023
$ws.Range("L61").Value = 2804.6667
$ws.Range("M61").Value = -811
$ws.Range("N61").Value = -3228.6667
$ws.Range("H63").Value = 4266.25
$ws.Range("I63").Value = 2573.3333
$ws.Range("J63").Value = 6442.857
$ws.Range("K63").Value = 2573.3333
$ws.Range("L63").Value = 6442.857
$ws.Range("M63").Value = -1887.3333
$ws.Range("N63").Value = -7814.857
$ws.Range("H66").Value = 4266.25
$ws.Range("I66").Value = 2573.3333
$ws.Range("J66").Value = 6442.857
$ws.Range("K66").Value = 12866.6665
$ws.Range("L66").Value = 32214.285
$ws.Range("M66").Value = -9434.666499999999
$ws.Range("N66").Value = -39078.285
$ws.Range("H132").Value = 13903971
$ws.Range("I132").Value = 19231746
$ws.Range("J132").Value = 51754.15
$ws.Range("K132").Value = 57695238
$ws.Range("L132").Value = 155262.45
$ws.Range("M132").Value = -57692708
$ws.Range("N132").Value = -160322.45
$ws.Range("H136").Value = 1720.174
$ws.Range("I136").Value = 1023
$ws.Range("J136").Value = 2804.6667
$ws.Range("K136").Value = 3069
$ws.Range("L136").Value = 8414.000100000001
$ws.Range("M136").Value = -519
$ws.Range("N136").Value = -13514.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2075.353
$ws.Range("I99").Value = 1460.909
$ws.Range("J99").Value = 3201.8333
$ws.Range("K99").Value = 1460.909
$ws.Range("L99").Value = 3201.8333
$ws.Range("M99").Value = 37.09099999999989
$ws.Range("N99").Value = -6197.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3379.422
$ws.Range("I31").Value = 4460.552
$ws.Range("J31").Value = 1419.875
$ws.Range("K31").Value = 4460.552
$ws.Range("L31").Value = 1419.875
$ws.Range("M31").Value = -4165.552
$ws.Range("N31").Value = -2009.875
$ws.Range("H34").Value = 3379.422
$ws.Range("I34").Value = 4460.552
$ws.Range("J34").Value = 1419.875
$ws.Range("K34").Value = 4460.552
$ws.Range("L34").Value = 1419.875
$ws.Range("M34").Value = -4258.552
$ws.Range("N34").Value = -1823.875
$ws.Range("H99").Value = 294577.44
$ws.Range("I99").Value = 437277.75
$ws.Range("J99").Value = 1666.3158
$ws.Range("K99").Value = 437277.75
$ws.Range("L99").Value = 1666.3158
$ws.Range("M99").Value = -435779.75
$ws.Range("N99").Value = -4662.3158
$ws.Range("H126").Value = 294577.44
$ws.Range("I126").Value = 437277.75
$ws.Range("J126").Value = 1666.3158
$ws.Range("K126").Value = 1311833.25
$ws.Range("L126").Value = 4998.9474
$ws.Range("M126").Value = -1309363.25
$ws.Range("N126").Value = -9938.947400000001
$ws.Range("H132").Value = 56061.105
$ws.Range("I132").Value = 2513.5
$ws.Range("K132").Value = 7540.5
$ws.Range("M132").Value = -5010.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 20095546
$ws.Range("I70").Value = 31254734
$ws.Range("J70").Value = 9009.9
$ws.Range("K70").Value = 31254734
$ws.Range("L70").Value = 9009.9
$ws.Range("M70").Value = -31254464
$ws.Range("N70").Value = -9549.9
$ws.Range("H73").Value = 20095546
$ws.Range("I73").Value = 31254734
$ws.Range("J73").Value = 9009.9
$ws.Range("K73").Value = 31254734
$ws.Range("L73").Value = 9009.9
$ws.Range("M73").Value = -31253798
$ws.Range("N73").Value = -10881.9
$ws.Range("H80").Value = 114359
$ws.Range("J80").Value = 3825.625
$ws.Range("L80").Value = 3825.625
$ws.Range("N80").Value = -5821.625
$ws.Range("H83").Value = 114359
$ws.Range("J83").Value = 3825.625
$ws.Range("L83").Value = 19128.125
$ws.Range("N83").Value = -29112.125
$ws.Range("H102").Value = 1378.9697
$ws.Range("I102").Value = 1066.4615
$ws.Range("J102").Value = 2539.7144
$ws.Range("K102").Value = 1066.4615
$ws.Range("L102").Value = 2539.7144
$ws.Range("M102").Value = 555.5385000000001
$ws.Range("N102").Value = -5783.7144
$ws.Range("H132").Value = 40795.46
$ws.Range("I132").Value = 1672.5
$ws.Range("J132").Value = 74329.42999999999
$ws.Range("K132").Value = 5017.5
$ws.Range("L132").Value = 222988.29
$ws.Range("M132").Value = -2487.5
$ws.Range("N132").Value = -228048.29

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2245.75
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("H126").Value = 2245.75
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N7").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1748.4706
$ws.Range("I126").Value = 1325
$ws.Range("J126").Value = 3124.75
$ws.Range("K126").Value = 3975
$ws.Range("L126").Value = 9374.25
$ws.Range("M126").Value = -1505
$ws.Range("N126").Value = -14314.25
$ws.Range("H132").Value = 69464020
$ws.Range("I132").Value = 113000984
$ws.Range("J132").Value = 2484059
$ws.Range("K132").Value = 339002952
$ws.Range("L132").Value = 7452177
$ws.Range("M132").Value = -339000422
$ws.Range("N132").Value = -7457237

Write-Output "Applied all cell updates."
